# updated CB_API and Dash
# Refresh the "future_occ" occupancy values for Feb/Mar/Apr/May 2024 rows
# (rows 2-5) to reflect the latest CB_API and Dash pull.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (2024-02)
$ws.Range("W2").Value = 0
$ws.Range("X2").Value = 0.45
$ws.Range("Y2").Value = 0.4666666666666667
$ws.Range("Z2").Value = 0.6166666666666667
$ws.Range("AA2").Value = 0.5833333333333334

# Row 3 (2024-03)
$ws.Range("C3").Value = 0.7833333333333333
$ws.Range("D3").Value = 0.85
$ws.Range("E3").Value = 0.5333333333333333
$ws.Range("M3").Value = 0.8833333333333333
$ws.Range("N3").Value = 0.8666666666666667
$ws.Range("O3").Value = 0.85
$ws.Range("P3").Value = 0.8333333333333334
$ws.Range("Q3").Value = 0.7833333333333333
$ws.Range("R3").Value = 0.55
$ws.Range("X3").Value = 0.3
$ws.Range("Y3").Value = 0.1833333333333333
$ws.Range("Z3").Value = 0.2166666666666667
$ws.Range("AA3").Value = 0.3666666666666666
$ws.Range("AC3").Value = 0.2833333333333333
$ws.Range("AD3").Value = 0.3
$ws.Range("AE3").Value = 0.3833333333333334
$ws.Range("AF3").Value = 0.4833333333333333
$ws.Range("AG3").Value = 0.2333333333333333

# Row 4 (2024-04)
$ws.Range("C4").Value = 0.4333333333333333
$ws.Range("E4").Value = 0.7166666666666667
$ws.Range("G4").Value = 0.65
$ws.Range("W4").Value = 0.1333333333333333
$ws.Range("X4").Value = 0.1166666666666667

# Row 5 (2024-05)
$ws.Range("S5").Value = 0.08333333333333333
$ws.Range("T5").Value = 0.1
$ws.Range("U5").Value = 0.05
